# Horarios Linea 141 - actualizacion 05:18:23 (commit: "Horarios actualizados Linea 141 - 421")
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173):
#   - refresh the "Ultima actualizacion" / "Total filas" header rows
#   - insert one new scrape row into the LP1912 sheet's existing table (at row 21)
#   - append the newest scrape rows to the bottom of each table

$wb = $excel.ActiveWorkbook

function Set-HeaderRows {
    param($ws, [string]$updatedAt, [int]$totalRows)
    $ws.Range("A2").Value = "Última actualización: " + $updatedAt
    $ws.Range("A3").Value = "Total filas: " + $totalRows
}

function Write-DataRow {
    param($ws, $rowNum, $horaScrap, $horaLlegada, $linea, $minutos, $parada)
    $ws.Cells.Item($rowNum, 1).Value = $horaScrap
    $ws.Cells.Item($rowNum, 2).Value = $horaLlegada
    $ws.Cells.Item($rowNum, 3).Value = $linea
    $ws.Cells.Item($rowNum, 4).Value = $minutos
    $ws.Cells.Item($rowNum, 5).Value = $parada
}

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

Set-HeaderRows $ws1 "05:18:23" 39

# A new scrape row belongs right after the existing "05:17" arrival (row 20),
# ahead of the "05:22" arrival that used to be row 21 - so push rows 21..37
# down one slot and drop the new row into the gap left at row 21.
$ws1.Rows.Item(21).Insert()
Write-DataRow $ws1 21 "05:18:23" "05:20" "14_ABASTO" 2 "LP1912"

# Newest scrape also appends six brand-new arrivals onto the end of the table
# (previously rows 6..38, now extending through row 44).
$sheet1NewRows = @(
    @(39, "05:18:23", "06:58", "10_OLMOS", 100, "LP1912"),
    @(40, "05:18:23", "06:59", "14_ABASTO", 101, "LP1912"),
    @(41, "05:18:23", "07:05", "15_ABASTO", 107, "LP1912"),
    @(42, "05:18:23", "07:07", "225_GOMEZ", 109, "LP1912"),
    @(43, "05:18:23", "07:11", "215A_EL PATO", 113, "LP1912"),
    @(44, "05:18:23", "07:15", "11_ETCHEVERRY", 117, "LP1912")
)
foreach ($row in $sheet1NewRows) {
    Write-DataRow $ws1 $row[0] $row[1] $row[2] $row[3] $row[4] $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

Set-HeaderRows $ws2 "05:18:23" 11
Write-DataRow $ws2 16 "05:18:23" "07:11" "215A_EL PATO" 113 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

Set-HeaderRows $ws3 "05:18:23" 8
Write-DataRow $ws3 13 "05:18:23" "07:00" "215B_LP-P MOR-1 Y 57" 102 "L6173"
